$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization rows appended to the bottom of the sheet.
$rows = @(
    @("ItemData.Item.1000006", "솔라리"),
    @("MapNpcData.MapNpcMenu.1000000.1", "컷신1 재성"),
    @("MapNpcData.MapNpcMenu.1000000.2", "컷신2 재생")
)

$startRow = 11
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
}
